# DPLKINV001-008 - Setup PIC Broker Investasi: replace the long step-by-step
# "Tambah/View/Ubah/Hapus" instructions in column D with short scenario
# labels (mirrors the SCENARIO_DESC convention used elsewhere in the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "Tambah Setup PIC Broker"
$ws.Range("D3").Value = "View Setup PIC Broker"
$ws.Range("D4").Value = "Ubah Setup PIC Broker"
$ws.Range("D5").Value = "Hapus Setup PIC Broker"

# The D column keeps its wrap-text style, but the new single-line labels no
# longer need the tall multi-line row heights - shrink rows 2,3,5 to 30pt
# and let row 4 fall back to the sheet's default height.
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).RowHeight = 30

# Saved view state moves the active selection to D5.
$ws.Range("D5").Select()
